$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header value updates
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: delete B2 and D2 entirely, update C2 and E2
$ws.Range("B2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = -4.5192477786255836
$ws.Range("E2").Value = -1.8852598986136755

# Row 3 updates
$ws.Range("B3").Value = -5.4378866419480616
$ws.Range("C3").Value = 3.2813695555772853
$ws.Range("D3").Value = -2.7961996361685308
$ws.Range("E3").Value = 10.479197686557249

# Update selection
$ws.Range("B1:E3").Select()
